# Add the journal citation after "Demo" in the title line.
#
# The title paragraph reads:
#   "ReadMe document for source code and Demo"
# and needs to become:
#   "ReadMe document for source code and Demo (Rogers, Al Husseini et al., Circ Res 2021)"
#
# "Demo" also appears many more times later in the document (e.g. "Demo
# folder", "Demo_Code.ipynb", "Demo Jupyter file", ...), so the search is
# scoped to just the first paragraph (the document title) to make sure only
# that one occurrence is touched.

$d = $word.ActiveDocument

$titlePara = $d.Paragraphs(1)
$searchRange = $titlePara.Range

$found = $searchRange.Find.Execute("Demo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the found "Demo" and insert the new text right
    # after it, as its own run, leaving the existing "Demo" run untouched.
    $searchRange.Collapse(0)
    $searchRange.InsertAfter(" (Rogers, Al Husseini et al., Circ Res 2021)")
}
